$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 3) with results for a ticker that now has a verdict of "Buy".
$ws.Cells.Item(3, 1).Value = 42632.883506944447
$ws.Cells.Item(3, 2).Value = 9
$ws.Cells.Item(3, 3).Value = "Buy"
$ws.Cells.Item(3, 4).Value = 24
$ws.Cells.Item(3, 5).Value = 15731
$ws.Cells.Item(3, 6).Value = 1727
$ws.Cells.Item(3, 7).Value = 61
$ws.Cells.Item(3, 8).Value = 36
$ws.Cells.Item(3, 9).Value = 83
$ws.Cells.Item(3, 10).Value = 16
$ws.Cells.Item(3, 11).Value = 9198
$ws.Cells.Item(3, 12).Value = 246
$ws.Cells.Item(3, 13).Value = 145
$ws.Cells.Item(3, 14).Value = 25
$ws.Cells.Item(3, 15).Value = 5
$ws.Cells.Item(3, 16).Value = "Bag"
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0.87
$ws.Cells.Item(3, 19).Value = 0.0351
$ws.Cells.Item(3, 19).NumberFormat = $ws.Cells.Item(2, 19).NumberFormat
$ws.Cells.Item(3, 20).Value = -2.08
$ws.Cells.Item(3, 21).Value = 15.16
$ws.Cells.Item(3, 22).Value = "N/A"
$ws.Cells.Item(3, 23).Value = 0
